$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3428.8462
$ws.Range("I6").Value = 399.875
$ws.Range("K6").Value = 1199.625
$ws.Range("M6").Value = -1087.625

$ws.Range("H42").Value = 174
$ws.Range("I42").Value = 42
$ws.Range("J42").Value = 284
$ws.Range("K42").Value = 126
$ws.Range("L42").Value = 852
$ws.Range("M42").Value = 104
$ws.Range("N42").Value = -1312

$ws.Range("H98").Value = 113185
$ws.Range("I98").Value = 144869.28
$ws.Range("J98").Value = 2290
$ws.Range("K98").Value = 144869.28
$ws.Range("L98").Value = 2290
$ws.Range("M98").Value = -143371.28
$ws.Range("N98").Value = -5286

$ws.Range("H103").Value = 46404.547
$ws.Range("I103").Value = 125825
$ws.Range("J103").Value = 1021.4286
$ws.Range("K103").Value = 377475
$ws.Range("L103").Value = 3064.2858
$ws.Range("M103").Value = -376889
$ws.Range("N103").Value = -4236.2858

$ws.Range("H122").Value = 113185
$ws.Range("I122").Value = 144869.28
$ws.Range("J122").Value = 2290
$ws.Range("K122").Value = 434607.84
$ws.Range("L122").Value = 6870
$ws.Range("M122").Value = -432157.84
$ws.Range("N122").Value = -11770

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 4514.375
$ws.Range("I25").Value = 2685.8333
$ws.Range("J25").Value = 10000
$ws.Range("K25").Value = 2685.8333
$ws.Range("L25").Value = 10000
$ws.Range("M25").Value = -2283.8333
$ws.Range("N25").Value = -10804

$ws.Range("H45").Value = 1278.7368
$ws.Range("I45").Value = 1050.1
$ws.Range("K45").Value = 1050.1
$ws.Range("M45").Value = -673.0999999999999

$ws.Range("H122").Value = 2786.75
$ws.Range("I122").Value = 3103
$ws.Range("J122").Value = 2470.5
$ws.Range("K122").Value = 9309
$ws.Range("L122").Value = 7411.5
$ws.Range("M122").Value = -6859
$ws.Range("N122").Value = -12311.5

$ws.Range("H132").Value = 2295.4062
$ws.Range("I132").Value = 1669.9445
$ws.Range("J132").Value = 3099.5715
$ws.Range("K132").Value = 5009.833500000001
$ws.Range("L132").Value = 9298.7145
$ws.Range("M132").Value = -2479.833500000001
$ws.Range("N132").Value = -14358.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws.Range("H134").Value = 1147.5682
$ws.Range("I134").Value = 1022.825
$ws.Range("J134").Value = 2395
$ws.Range("K134").Value = 3068.475
$ws.Range("L134").Value = 7185
$ws.Range("M134").Value = -533.4750000000004
$ws.Range("N134").Value = -12255

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 44666.668
$ws.Range("J18").Value = 44666.668
$ws.Range("L18").Value = 44666.668
$ws.Range("N18").Value = -45126.668

$ws.Range("H36").Value = 11877.875
$ws.Range("I36").Value = 6841
$ws.Range("J36").Value = 14900
$ws.Range("K36").Value = 6841
$ws.Range("L36").Value = 14900
$ws.Range("M36").Value = -6453
$ws.Range("N36").Value = -15676

$ws.Range("H40").Value = 11877.875
$ws.Range("I40").Value = 6841
$ws.Range("J40").Value = 14900
$ws.Range("K40").Value = 6841
$ws.Range("L40").Value = 14900
$ws.Range("M40").Value = -6681
$ws.Range("N40").Value = -15220

$ws.Range("H99").Value = 142858240
$ws.Range("I99").Value = 1359.4
$ws.Range("J99").Value = 500000450
$ws.Range("K99").Value = 1359.4
$ws.Range("L99").Value = 500000450
$ws.Range("M99").Value = 138.5999999999999
$ws.Range("N99").Value = -500003446

$ws.Range("H117").Value = 41000
$ws.Range("J117").Value = 41000
$ws.Range("L117").Value = 41000
$ws.Range("N117").Value = -50178

$ws.Range("H126").Value = 142858240
$ws.Range("I126").Value = 1359.4
$ws.Range("J126").Value = 500000450
$ws.Range("K126").Value = 4078.2
$ws.Range("L126").Value = 1500001350
$ws.Range("M126").Value = -1608.2
$ws.Range("N126").Value = -1500006290

$ws.Range("H132").Value = 3049.7693
$ws.Range("I132").Value = 2345.4
$ws.Range("K132").Value = 7036.200000000001
$ws.Range("M132").Value = -4506.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 16326
$ws.Range("J51").Value = 16326
$ws.Range("L51").Value = 16326
$ws.Range("N51").Value = -17344

$ws.Range("H102").Value = 2271.4614
$ws.Range("I102").Value = 1766.2593
$ws.Range("J102").Value = 3408.1667
$ws.Range("K102").Value = 1766.2593
$ws.Range("L102").Value = 3408.1667
$ws.Range("M102").Value = -144.2592999999999
$ws.Range("N102").Value = -6652.1667

$ws.Range("H108").Value = 33000
$ws.Range("J108").Value = 33000
$ws.Range("L108").Value = 33000
$ws.Range("N108").Value = -40680

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4451.3
$ws.Range("I7").Value = 4644
$ws.Range("J7").Value = 4001.6667
$ws.Range("K7").Value = 4644
$ws.Range("L7").Value = 4001.6667
$ws.Range("M7").Value = -4532
$ws.Range("N7").Value = -4225.6667

$ws.Range("H40").Value = 2413.5334
$ws.Range("I40").Value = 2165.3044
$ws.Range("J40").Value = 3229.1428
$ws.Range("K40").Value = 2165.3044
$ws.Range("L40").Value = 3229.1428
$ws.Range("M40").Value = -2029.3044
$ws.Range("N40").Value = -3501.1428

$ws.Range("H120").Value = 30000
$ws.Range("J120").Value = 30000
$ws.Range("L120").Value = 30000
$ws.Range("N120").Value = -39676

$ws.Range("H122").Value = 3156.7297
$ws.Range("I122").Value = 3041.5173
$ws.Range("J122").Value = 3574.375
$ws.Range("K122").Value = 9124.5519
$ws.Range("L122").Value = 10723.125
$ws.Range("M122").Value = -6674.5519
$ws.Range("N122").Value = -15623.125

$ws.Range("H126").Value = 4451.3
$ws.Range("I126").Value = 4644
$ws.Range("J126").Value = 4001.6667
$ws.Range("K126").Value = 13932
$ws.Range("L126").Value = 12005.0001
$ws.Range("M126").Value = -11462
$ws.Range("N126").Value = -16945.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 8080.2
$ws.Range("J45").Value = 8080.2
$ws.Range("L45").Value = 8080.2
$ws.Range("N45").Value = -9062.200000000001

$ws.Range("H132").Value = 2037.6666
$ws.Range("I132").Value = 1337.7142
$ws.Range("J132").Value = 2650.125
$ws.Range("K132").Value = 4013.1426
$ws.Range("L132").Value = 7950.375
$ws.Range("M132").Value = -1483.1426
$ws.Range("N132").Value = -13010.375
